$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CreatedAt date-text formatting for existing rows (shared string)
$ws.Range('F2').Value = '9/15/2025 14:54'
$ws.Range('F3').Value = '9/15/2025 14:54'
$ws.Range('F4').Value = '9/15/2025 14:54'

# New product rows (5-10): Name / Description / Price / QtyInitial / QtySold
$ws.Cells.Item(5,1).Value = 'Nước Cam Ép'
$ws.Cells.Item(5,2).Value = 'Nước cam tươi nguyên chất'
$ws.Cells.Item(5,3).Value = 15000
$ws.Cells.Item(5,4).Value = 100
$ws.Cells.Item(5,5).Value = 35

$ws.Cells.Item(6,1).Value = 'Nước Chanh Tươi'
$ws.Cells.Item(6,2).Value = 'Nước chanh mát lạnh giải khát'
$ws.Cells.Item(6,3).Value = 12000
$ws.Cells.Item(6,4).Value = 90
$ws.Cells.Item(6,5).Value = 28

$ws.Cells.Item(7,1).Value = 'Nước Dừa Xiêm'
$ws.Cells.Item(7,2).Value = 'Nước dừa xiêm ngọt thanh'
$ws.Cells.Item(7,3).Value = 20000
$ws.Cells.Item(7,4).Value = 70
$ws.Cells.Item(7,5).Value = 25

$ws.Cells.Item(8,1).Value = 'Sữa Tươi'
$ws.Cells.Item(8,2).Value = 'Sữa tươi tiệt trùng nguyên chất'
$ws.Cells.Item(8,3).Value = 18000
$ws.Cells.Item(8,4).Value = 80
$ws.Cells.Item(8,5).Value = 30

$ws.Cells.Item(9,1).Value = 'Sữa Đậu Nành'
$ws.Cells.Item(9,2).Value = 'Thức uống từ đậu nành bổ dưỡng'
$ws.Cells.Item(9,3).Value = 12000
$ws.Cells.Item(9,4).Value = 90
$ws.Cells.Item(9,5).Value = 22

$ws.Cells.Item(10,1).Value = 'Trà Sữa Trân Châu'
$ws.Cells.Item(10,2).Value = 'Trà sữa kèm trân châu dai ngon'
$ws.Cells.Item(10,3).Value = 35000
$ws.Cells.Item(10,4).Value = 100
$ws.Cells.Item(10,5).Value = 50

# New rows: CreatedAt left blank (explicit empty text)
$ws.Cells.Item(5,6).Value = ''
$ws.Cells.Item(6,6).Value = ''
$ws.Cells.Item(7,6).Value = ''
$ws.Cells.Item(8,6).Value = ''
$ws.Cells.Item(9,6).Value = ''
$ws.Cells.Item(10,6).Value = ''

# New rows: image URL columns (G/H/I) + hyperlinks
$ws.Cells.Item(5,7).Value = 'https://bavifoods.com/thumbs/740x740x1/upload/product/cam-ep-5018.jpg'
$ws.Hyperlinks.Add($ws.Cells.Item(5,7), 'https://bavifoods.com/thumbs/740x740x1/upload/product/cam-ep-5018.jpg', '', '', 'https://bavifoods.com/thumbs/740x740x1/upload/product/cam-ep-5018.jpg')
$ws.Cells.Item(5,8).Value = 'https://bavifoods.com/thumbs/740x740x1/upload/product/cam-ep-5018.jpg'
$ws.Hyperlinks.Add($ws.Cells.Item(5,8), 'https://bavifoods.com/thumbs/740x740x1/upload/product/cam-ep-5018.jpg', '', '', 'https://bavifoods.com/thumbs/740x740x1/upload/product/cam-ep-5018.jpg')
$ws.Cells.Item(5,9).Value = 'https://cdn-i.vtcnews.vn/files/news/2019/01/22/-145625.jpg'
$ws.Hyperlinks.Add($ws.Cells.Item(5,9), 'https://cdn-i.vtcnews.vn/files/news/2019/01/22/-145625.jpg', '', '', 'https://cdn-i.vtcnews.vn/files/news/2019/01/22/-145625.jpg')

$ws.Cells.Item(6,7).Value = 'https://media.baobinhphuoc.com.vn/upload/news/5_2023/img_8476_06413001052023.jpeg'
$ws.Hyperlinks.Add($ws.Cells.Item(6,7), 'https://media.baobinhphuoc.com.vn/upload/news/5_2023/img_8476_06413001052023.jpeg', '', '', 'https://media.baobinhphuoc.com.vn/upload/news/5_2023/img_8476_06413001052023.jpeg')
$ws.Cells.Item(6,8).Value = 'https://media.vov.vn/sites/default/files/styles/large/public/2023-06/nuoc_chanh_5.jpg'
$ws.Hyperlinks.Add($ws.Cells.Item(6,8), 'https://media.vov.vn/sites/default/files/styles/large/public/2023-06/nuoc_chanh_5.jpg', '', '', 'https://media.vov.vn/sites/default/files/styles/large/public/2023-06/nuoc_chanh_5.jpg')
$ws.Cells.Item(6,9).Value = 'https://suckhoedoisong.qltns.mediacdn.vn/324455921873985536/2022/4/18/uong-nuoc-chanh-moi-ngay-co-tot-khong-va-uong-khi-nao-chanh1-1592466583-666-width1024height768-16502694973802093569436.jpg'
$ws.Hyperlinks.Add($ws.Cells.Item(6,9), 'https://suckhoedoisong.qltns.mediacdn.vn/324455921873985536/2022/4/18/uong-nuoc-chanh-moi-ngay-co-tot-khong-va-uong-khi-nao-chanh1-1592466583-666-width1024height768-16502694973802093569436.jpg', '', '', 'https://suckhoedoisong.qltns.mediacdn.vn/324455921873985536/2022/4/18/uong-nuoc-chanh-moi-ngay-co-tot-khong-va-uong-khi-nao-chanh1-1592466583-666-width1024height768-16502694973802093569436.jpg')

$ws.Cells.Item(7,7).Value = 'Chi-em-thi-nhau-lung-mua-dua-xiem-ve-uong-sau-tiem-phong-cua-hang-moi-ngay-ban-5000-qua-1-1631524190-570-width650height431.jpg (650×431)'
$ws.Hyperlinks.Add($ws.Cells.Item(7,7), 'Chi-em-thi-nhau-lung-mua-dua-xiem-ve-uong-sau-tiem-phong-cua-hang-moi-ngay-ban-5000-qua-1-1631524190-570-width650height431.jpg (650×431)', '', '', 'Chi-em-thi-nhau-lung-mua-dua-xiem-ve-uong-sau-tiem-phong-cua-hang-moi-ngay-ban-5000-qua-1-1631524190-570-width650height431.jpg (650×431)')
$ws.Cells.Item(7,8).Value = 'medium_20200513_094458_574364_nuoc_dua_max_1800x1800_jpg_095dc5e7ad.jpg (750×563)'
$ws.Hyperlinks.Add($ws.Cells.Item(7,8), 'medium_20200513_094458_574364_nuoc_dua_max_1800x1800_jpg_095dc5e7ad.jpg (750×563)', '', '', 'medium_20200513_094458_574364_nuoc_dua_max_1800x1800_jpg_095dc5e7ad.jpg (750×563)')
$ws.Cells.Item(7,9).Value = 'coconut-water-benefits-17218412875751213756362.jpg (800×562)'
$ws.Hyperlinks.Add($ws.Cells.Item(7,9), 'coconut-water-benefits-17218412875751213756362.jpg (800×562)', '', '', 'coconut-water-benefits-17218412875751213756362.jpg (800×562)')

$ws.Cells.Item(8,7).Value = 'https://file.hstatic.net/1000199715/file/uong-sua-sau-sinh-1_90f6b928e6084e7e87c4e7a89e1b1be3_grande.jpg'
$ws.Hyperlinks.Add($ws.Cells.Item(8,7), 'https://file.hstatic.net/1000199715/file/uong-sua-sau-sinh-1_90f6b928e6084e7e87c4e7a89e1b1be3_grande.jpg', '', '', 'https://file.hstatic.net/1000199715/file/uong-sua-sau-sinh-1_90f6b928e6084e7e87c4e7a89e1b1be3_grande.jpg')
$ws.Cells.Item(8,8).Value = 'khi_nao_nen_cho_be_uong_sua_1_4401cf044a.jpg (800×600)'
$ws.Hyperlinks.Add($ws.Cells.Item(8,8), 'khi_nao_nen_cho_be_uong_sua_1_4401cf044a.jpg (800×600)', '', '', 'khi_nao_nen_cho_be_uong_sua_1_4401cf044a.jpg (800×600)')

$ws.Cells.Item(9,7).Value = 'https://suckhoedoisong.qltns.mediacdn.vn/324455921873985536/2025/2/21/dau-nanh-1-1740125251401155246723.jpg'
$ws.Hyperlinks.Add($ws.Cells.Item(9,7), 'https://suckhoedoisong.qltns.mediacdn.vn/324455921873985536/2025/2/21/dau-nanh-1-1740125251401155246723.jpg', '', '', 'https://suckhoedoisong.qltns.mediacdn.vn/324455921873985536/2025/2/21/dau-nanh-1-1740125251401155246723.jpg')
$ws.Cells.Item(9,8).Value = 'glass-soy-milk_20dc83bb32164c49bd11a7d7b60b717b_grande.jpg (600×377)'
$ws.Hyperlinks.Add($ws.Cells.Item(9,8), 'glass-soy-milk_20dc83bb32164c49bd11a7d7b60b717b_grande.jpg (600×377)', '', '', 'glass-soy-milk_20dc83bb32164c49bd11a7d7b60b717b_grande.jpg (600×377)')
$ws.Cells.Item(9,9).Value = 'may-lam-sua-dau-nanh-1-1412734006024.jpg (500×455)'
$ws.Hyperlinks.Add($ws.Cells.Item(9,9), 'may-lam-sua-dau-nanh-1-1412734006024.jpg (500×455)', '', '', 'may-lam-sua-dau-nanh-1-1412734006024.jpg (500×455)')

$ws.Cells.Item(10,7).Value = 'https://baothainguyen.vn/file/e7837c027f6ecd14017ffa4e5f2a0e34/032023/1-boba-tea-recipe-using-fresh-tapioca-pearls-1024x1024-1677809524112848165864_20230305161118.jpeg'
$ws.Hyperlinks.Add($ws.Cells.Item(10,7), 'https://baothainguyen.vn/file/e7837c027f6ecd14017ffa4e5f2a0e34/032023/1-boba-tea-recipe-using-fresh-tapioca-pearls-1024x1024-1677809524112848165864_20230305161118.jpeg', '', '', 'https://baothainguyen.vn/file/e7837c027f6ecd14017ffa4e5f2a0e34/032023/1-boba-tea-recipe-using-fresh-tapioca-pearls-1024x1024-1677809524112848165864_20230305161118.jpeg')
$ws.Cells.Item(10,8).Value = 'https://www.cet.edu.vn/wp-content/uploads/2018/04/tra-sua-tu-lam.jpg'
$ws.Hyperlinks.Add($ws.Cells.Item(10,8), 'https://www.cet.edu.vn/wp-content/uploads/2018/04/tra-sua-tu-lam.jpg', '', '', 'https://www.cet.edu.vn/wp-content/uploads/2018/04/tra-sua-tu-lam.jpg')

# Formatting: wrap text across the data block (rows 2-10)
$dataBlock = $ws.Range('A2:I10')
$dataBlock.WrapText = $true

# Right-align numeric columns C:E for rows 2-10
$numBlock = $ws.Range('C2:E10')
$numBlock.HorizontalAlignment = -4152

# CreatedAt column formatting
# Existing rows 2-4: '@' text format, right aligned
$ws.Range('F2:F4').NumberFormat = '@'
$ws.Range('F2:F4').HorizontalAlignment = -4152
$ws.Range('F2:F4').WrapText = $true
# New rows 5-10: general format, left aligned, no wrap
$ws.Range('F5:F10').HorizontalAlignment = -4131
$ws.Range('F5:F10').WrapText = $false

# I8 / I10 stay empty but share the general/no-wrap style used for F5:F10
$ws.Range('I8').WrapText = $false
$ws.Range('I8').HorizontalAlignment = -4131
$ws.Range('I10').WrapText = $false
$ws.Range('I10').HorizontalAlignment = -4131

# G5: white font color (URL hidden against white background)
$ws.Range('G5').Font.Color = 16777215

# Row heights
$ws.Rows.Item(2).RowHeight = 35.05
$ws.Rows.Item(3).RowHeight = 35.05
$ws.Rows.Item(4).RowHeight = 35.05
$ws.Rows.Item(5).RowHeight = 46.25
$ws.Rows.Item(6).RowHeight = 46.25
$ws.Rows.Item(7).RowHeight = 35.05
$ws.Rows.Item(8).RowHeight = 15
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(10).RowHeight = 15

# Column A width
$ws.Columns.Item(1).ColumnWidth = 29.7

# Active cell selection
$ws.Range('G4').Select()

Write-Host 'done'
